$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(2, 2).Range.Text = "227"
$tbl.Cell(2, 3).Range.Text = "138"

$tbl.Cell(3, 2).Range.Text = "199 (97,5)"
$tbl.Cell(3, 3).Range.Text = "122 (94,6)"
$tbl.Cell(3, 4).Range.Text = "0,226"

$tbl.Cell(4, 2).Range.Text = "4 ( 2,0)"
$tbl.Cell(4, 3).Range.Text = "1 ( 0,8)"
$tbl.Cell(4, 4).Range.Text = "0,652"

$tbl.Cell(5, 2).Range.Text = "197 (96,6)"
$tbl.Cell(5, 3).Range.Text = "122 (93,8)"
$tbl.Cell(5, 4).Range.Text = "0,283"

$tbl.Cell(6, 1).Range.Text = "Rivaroxabana (%)"
$tbl.Cell(6, 2).Range.Text = "44 (21,5)"
$tbl.Cell(6, 3).Range.Text = "28 (21,5)"
$tbl.Cell(6, 4).Range.Text = "1,000"

$tbl.Cell(7, 2).Range.Text = "121 (63,4)"
$tbl.Cell(7, 3).Range.Text = "56 (47,9)"
$tbl.Cell(7, 4).Range.Text = "0,009"
